$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109; this pushes the existing rows 109-150
# down to 110-151 (preserving all of their data/formatting), matching the
# dimension growing from A1:T150 to A1:T151.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with its data (same constant columns
# as the rest of the "Pi\u00f1a" / "Caramelo" / "Ecuador" block, plus the
# new record's own date, quality, volume, price and unit values).
$ws.Range("A109").Value = 5
$ws.Range("B109").Value = "Macroferia Regional de Talca"
$ws.Range("C109").Value = "Maule"
$ws.Range("D109").Value = 44468
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = "Tropicales y subtropicales"
$ws.Range("I109").Value = 100108005
$ws.Range("J109").Value = "Piña"
$ws.Range("K109").Value = "Caramelo"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 120
$ws.Range("N109").Value = 19000
$ws.Range("O109").Value = 19000
$ws.Range("P109").Value = 19000
$ws.Range("Q109").Value = "`$/caja 14 unidades"
$ws.Range("R109").Value = "Ecuador"
$ws.Range("S109").Value = 1357
$ws.Range("T109").Value = 14
